$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'303.21"
$ws.Range("E2").Value = "'4.97%"
$ws.Range("D3").Value = "'34.84"
$ws.Range("E3").Value = "'12.24%"
$ws.Range("D4").Value = "'5.141"
$ws.Range("E4").Value = "'4.27%"
$ws.Range("D5").Value = "'0.07756"
$ws.Range("E5").Value = "'5.04%"
$ws.Range("D6").Value = "'2.355"
$ws.Range("E6").Value = "'6.17%"
$ws.Range("D7").Value = "'8.019"
$ws.Range("E7").Value = "'4.23%"
$ws.Range("D8").Value = "'3.946"
$ws.Range("E8").Value = "'5.73%"
$ws.Range("D9").Value = "'0.9284"
$ws.Range("E9").Value = "'2.13%"
$ws.Range("D10").Value = "'0.1013"
$ws.Range("E10").Value = "'16.23%"
$ws.Range("D11").Value = "'0.1797"
$ws.Range("E11").Value = "'6.48%"
$ws.Range("D12").Value = "'0.08511"
$ws.Range("E12").Value = "'3.69%"
$ws.Range("E13").Value = "'6.55%"
$ws.Range("D14").Value = "'0.09894"
$ws.Range("E14").Value = "'-0.56%"
$ws.Range("D15").Value = "'0.001499"
$ws.Range("E15").Value = "'-0.09%"
$ws.Range("D16").Value = "'0.005760"
$ws.Range("E16").Value = "'-0.65%"
$ws.Range("E17").Value = "'-0.65%"
$ws.Range("E18").Value = "'4.52%"
$ws.Range("E19").Value = "'1.16%"
$ws.Range("E20").Value = "'1.05%"
$ws.Range("D21").Value = "'4.331"
$ws.Range("E21").Value = "'13.15%"
$ws.Range("D22").Value = "'0.2387"
$ws.Range("E22").Value = "'9.02%"
$ws.Range("E23").Value = "'0.31%"
$ws.Range("D24").Value = "'0.001216"
$ws.Range("E24").Value = "'0.46%"
$ws.Range("D25").Value = "'0.004465"
$ws.Range("E25").Value = "'7.75%"
$ws.Range("E26").Value = "'-0.18%"
$ws.Range("E27").Value = "'-0.19%"
$ws.Range("D39").Value = "'0.01776"
$ws.Range("E39").Value = "'12.26%"
$ws.Range("D40").Value = "'0.04746"
$ws.Range("E40").Value = "'6.28%"
$ws.Range("D41").Value = "'0.007747"
$ws.Range("E41").Value = "'5.38%"
$ws.Range("D42").Value = "'0.1412"
$ws.Range("E42").Value = "'6.59%"
$ws.Range("D43").Value = "'0.007092"
$ws.Range("E43").Value = "'-25.84%"
$ws.Range("D44").Value = "'0.002150"
$ws.Range("E44").Value = "'0.62%"
$ws.Range("D45").Value = "'0.009186"
$ws.Range("E45").Value = "'10.61%"
$ws.Range("D46").Value = "'0.00006120"
$ws.Range("E46").Value = "'-0.03%"
$ws.Range("E47").Value = "'-0.16%"
$ws.Range("D48").Value = "'2.726"
$ws.Range("E48").Value = "'29.41%"
$ws.Range("D49").Value = "'0.002000"
$ws.Range("E49").Value = "'-0.14%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.16%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.16%"
